# Regenerate ranking content per analysis module.
#
# The underlying UF rankings are recomputed independently for each module
# (qtd, tot_arrecadado, avg_arrecadado, max_arrecadado, txsucesso). For rows
# that are tied on the ranked value, the tie-break order can change between
# runs. Two modules are affected by this re-generation:
#
#   - "qtd": BA/PB are tied at 30 (rows 10-11); MA/RO/AM are tied at 4
#     (rows 20-22).
#   - "tx-sucesso": MA/XX/MT/RO are tied at 100 (rows 2-5).
#
# The B-column (value) data is untouched; only the state labels in column A
# are reordered within each tied group.

$wb = $excel.ActiveWorkbook

$wsQtd = $wb.Worksheets.Item("qtd")
$wsQtd.Range("A10").Value = "PB"
$wsQtd.Range("A11").Value = "BA"
$wsQtd.Range("A20").Value = "RO"
$wsQtd.Range("A21").Value = "AM"
$wsQtd.Range("A22").Value = "MA"

$wsTx = $wb.Worksheets.Item("tx-sucesso")
$wsTx.Range("A2").Value = "RO"
$wsTx.Range("A3").Value = "MT"
$wsTx.Range("A4").Value = "MA"
$wsTx.Range("A5").Value = "XX"
